$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 910953
$ws.Range("I2").Value2 = 1818556.5
$ws.Range("J2").Value2 = 3349.5
$ws.Range("K2").Value2 = 1818556.5
$ws.Range("L2").Value2 = 3349.5
$ws.Range("M2").Value2 = -1818443.5
$ws.Range("N2").Value2 = -3575.5

# Sheet ALC, Row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 34690
$ws.Range("J17").Value2 = 51035
$ws.Range("L17").Value2 = 153105
$ws.Range("N17").Value2 = -153441

# Sheet ALC, Row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 11277.096
$ws.Range("I62").Value2 = 13003.223
$ws.Range("J62").Value2 = 9982.5
$ws.Range("K62").Value2 = 13003.223
$ws.Range("L62").Value2 = 9982.5
$ws.Range("M62").Value2 = -12379.223
$ws.Range("N62").Value2 = -11230.5

# Sheet ALC, Row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value2 = 11277.096
$ws.Range("I65").Value2 = 13003.223
$ws.Range("J65").Value2 = 9982.5
$ws.Range("K65").Value2 = 65016.115
$ws.Range("L65").Value2 = 49912.5
$ws.Range("M65").Value2 = -61896.115
$ws.Range("N65").Value2 = -56152.5

# Sheet ALC, Row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value2 = 3733.2222
$ws.Range("J74").Value2 = 4000
$ws.Range("L74").Value2 = 4000
$ws.Range("N74").Value2 = -5872

# Sheet ALC, Row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value2 = 3733.2222
$ws.Range("J77").Value2 = 4000
$ws.Range("L77").Value2 = 20000
$ws.Range("N77").Value2 = -29360

# Sheet ALC, Row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value2 = 48281.81
$ws.Range("J92").Value2 = 899.2222
$ws.Range("L92").Value2 = 899.2222
$ws.Range("N92").Value2 = -3395.2222

# Sheet ALC, Row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value2 = 1534.2
$ws.Range("J96").Value2 = 2482.5
$ws.Range("L96").Value2 = 7447.5
$ws.Range("N96").Value2 = -10193.5

# Sheet ALC, Row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value2 = 1484.75
$ws.Range("I127").Value2 = 1484.75
$ws.Range("J127").Value2 = 0
$ws.Range("K127").Value2 = 4454.25
$ws.Range("L127").Value2 = 0
$ws.Range("M127").Value2 = 505.75
$ws.Range("N127").Value2 = ""

# Sheet ARM, Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 3314.8
$ws.Range("I45").Value2 = 2335.7856
$ws.Range("K45").Value2 = 2335.7856
$ws.Range("M45").Value2 = -1958.7856

# Sheet ARM, Row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value2 = 2165.6667
$ws.Range("I102").Value2 = 2165.6667
$ws.Range("K102").Value2 = 2165.6667
$ws.Range("M102").Value2 = -543.6667000000002

# Sheet ARM, Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value2 = 3042.4546
$ws.Range("I122").Value2 = 2746.7
$ws.Range("K122").Value2 = 8240.099999999999
$ws.Range("M122").Value2 = -5790.099999999999

# Sheet BSM, Row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value2 = 2699.1765
$ws.Range("I94").Value2 = 2453.6365
$ws.Range("J94").Value2 = 3149.3333
$ws.Range("K94").Value2 = 2453.6365
$ws.Range("L94").Value2 = 3149.3333
$ws.Range("M94").Value2 = -2002.6365
$ws.Range("N94").Value2 = -4051.3333

# Sheet BSM, Row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 3468
$ws.Range("I134").Value2 = 3179.5454
$ws.Range("J134").Value2 = 4525.6665
$ws.Range("K134").Value2 = 9538.636200000001
$ws.Range("L134").Value2 = 13576.9995
$ws.Range("M134").Value2 = -7003.636200000001
$ws.Range("N134").Value2 = -18646.9995

# Sheet CRP, Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3365

# Sheet CRP, Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 3365

# Sheet CRP, Row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value2 = 3210.2727
$ws.Range("I94").Value2 = 3121.25
$ws.Range("J94").Value2 = 3261.1428
$ws.Range("K94").Value2 = 3121.25
$ws.Range("L94").Value2 = 3261.1428
$ws.Range("M94").Value2 = -2670.25
$ws.Range("N94").Value2 = -4163.1428

# Sheet CRP, Row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value2 = 2721.8333
$ws.Range("I99").Value2 = 1350
$ws.Range("K99").Value2 = 1350
$ws.Range("M99").Value2 = 148

# Sheet CRP, Row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 500.25
$ws.Range("I107").Value2 = 480.29413
$ws.Range("K107").Value2 = 480.29413
$ws.Range("M107").Value2 = 1439.70587

# Sheet CRP, Row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value2 = 2411.6667
$ws.Range("I122").Value2 = 2411.6667
$ws.Range("K122").Value2 = 7235.000100000001
$ws.Range("M122").Value2 = -4785.000100000001

# Sheet CRP, Row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value2 = 2721.8333
$ws.Range("I126").Value2 = 1350
$ws.Range("K126").Value2 = 4050
$ws.Range("M126").Value2 = -1580

# Sheet CRP, Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 4680.75
$ws.Range("I132").Value2 = 4608
$ws.Range("K132").Value2 = 13824
$ws.Range("M132").Value2 = -11294

# Sheet CUL, Row 128
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value2 = 122197.4
$ws.Range("I128").Value2 = 122197.4
$ws.Range("K128").Value2 = 366592.2
$ws.Range("M128").Value2 = -361612.2

# Sheet CUL, Row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value2 = 1603.0834
$ws.Range("I139").Value2 = 1294.2727
$ws.Range("K139").Value2 = 3882.8181
$ws.Range("M139").Value2 = 1257.1819

# Sheet GSM, Row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value2 = 51481
$ws.Range("J15").Value2 = 51481
$ws.Range("L15").Value2 = 51481
$ws.Range("N15").Value2 = -52057

# Sheet GSM, Row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value2 = 51481
$ws.Range("J81").Value2 = 51481
$ws.Range("L81").Value2 = 51481
$ws.Range("N81").Value2 = -53477

# Sheet GSM, Row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value2 = 51481
$ws.Range("J84").Value2 = 51481
$ws.Range("L84").Value2 = 154443
$ws.Range("N84").Value2 = -164427

# Sheet GSM, Row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 1767.697
$ws.Range("I102").Value2 = 1344.5
$ws.Range("K102").Value2 = 1344.5
$ws.Range("M102").Value2 = 277.5

# Sheet GSM, Row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 2491.3
$ws.Range("I122").Value2 = 2058.6428
$ws.Range("J122").Value2 = 3500.8333
$ws.Range("K122").Value2 = 6175.928400000001
$ws.Range("L122").Value2 = 10502.4999
$ws.Range("M122").Value2 = -3725.928400000001
$ws.Range("N122").Value2 = -15402.4999

# Sheet GSM, Row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 5692.9
$ws.Range("I126").Value2 = 3567
$ws.Range("K126").Value2 = 10701
$ws.Range("M126").Value2 = -8231

# Sheet GSM, Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 53034.45
$ws.Range("I132").Value2 = 58549.39
$ws.Range("K132").Value2 = 175648.17
$ws.Range("M132").Value2 = -173118.17

# Sheet LTW, Row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 5671.4
$ws.Range("J7").Value2 = 4670
$ws.Range("L7").Value2 = 4670
$ws.Range("N7").Value2 = -4894

# Sheet LTW, Row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 409.33334
$ws.Range("I16").Value2 = 409.33334
$ws.Range("K16").Value2 = 409.33334
$ws.Range("M16").Value2 = -239.33334

# Sheet LTW, Row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 4838.125
$ws.Range("I40").Value2 = 4357.857
$ws.Range("K40").Value2 = 4357.857
$ws.Range("M40").Value2 = -4221.857

# Sheet LTW, Row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 2822.9443
$ws.Range("I61").Value2 = 2054.2
$ws.Range("K61").Value2 = 2054.2
$ws.Range("M61").Value2 = -1852.2

# Sheet LTW, Row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 1931.0344
$ws.Range("I93").Value2 = 1690.091
$ws.Range("J93").Value2 = 2078.2778
$ws.Range("K93").Value2 = 1690.091
$ws.Range("L93").Value2 = 2078.2778
$ws.Range("M93").Value2 = -442.0909999999999
$ws.Range("N93").Value2 = -4574.2778

# Sheet LTW, Row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value2 = 2822.9443
$ws.Range("I113").Value2 = 2054.2
$ws.Range("K113").Value2 = 2054.2
$ws.Range("M113").Value2 = 115.8000000000002

# Sheet LTW, Row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 3728.8572
$ws.Range("I122").Value2 = 3128.0908
$ws.Range("J122").Value2 = 4745.5386
$ws.Range("K122").Value2 = 9384.2724
$ws.Range("L122").Value2 = 14236.6158
$ws.Range("M122").Value2 = -6934.2724
$ws.Range("N122").Value2 = -19136.6158

# Sheet LTW, Row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value2 = 5671.4
$ws.Range("J126").Value2 = 4670
$ws.Range("L126").Value2 = 14010
$ws.Range("N126").Value2 = -18950

# Sheet WVR, Row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 870.1
$ws.Range("I107").Value2 = 850.2857
$ws.Range("K107").Value2 = 2550.8571
$ws.Range("M107").Value2 = -630.8571000000002

# Sheet WVR, Row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 7684.6523
$ws.Range("I122").Value2 = 8476.210999999999
$ws.Range("J122").Value2 = 3924.75
$ws.Range("K122").Value2 = 25428.633
$ws.Range("L122").Value2 = 11774.25
$ws.Range("M122").Value2 = -22978.633
$ws.Range("N122").Value2 = -16674.25

# Sheet WVR, Row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 133917.05
$ws.Range("I126").Value2 = 171621.77
$ws.Range("K126").Value2 = 514865.3099999999
$ws.Range("M126").Value2 = -512395.3099999999

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 119645.586
$ws.Range("I132").Value2 = 134331.86
$ws.Range("K132").Value2 = 402995.58
$ws.Range("M132").Value2 = -400465.58

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value2 = 4045.0454
$ws.Range("I136").Value2 = 4009.0952
$ws.Range("K136").Value2 = 12027.2856
$ws.Range("M136").Value2 = -9477.285600000001
